$d = $word.ActiveDocument

# 1) "587-1" -> "588-1" (appears in both "Выдано: студенту группы 587-1"
#    and "студент гр. 587-1"); replace both occurrences document-wide.
$d.Content.Find.Execute("587-1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "588-1", 2)

# 2) Update the due date from 03.05. to 13.10. and drop the trailing
#    space right after "работы:" (" работы: 03.05." -> " работы:13.10.").
$d.Content.Find.Execute(" работы: 03.05.", $true, $false, $false, $false, `
                         $false, $true, 1, $false, " работы:13.10.", 2)

# 3) Remove the stray "_GoBack" bookmark left over from the last edit
#    position.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
